$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Item")

# Row 2 holds the per-column type declarations for the Proto sheet.
# SizeX (E2) and SizeY (F2) were mis-typed as "string" - correct them to "int".
$ws.Range("E2").Value = "int"
$ws.Range("F2").Value = "int"

# Reflect the author's active selection at save time (E2:F2, active cell E2).
$ws.Range("E2:F2").Select()
